# Actualización automática 2025-10-31 16:30:08
#
# Source-data corrections ripple through three related sheets:
#   1) "VENTAS POR GRUPO"     - per-client sales broken down by product group
#   2) "VENTA MENSUAL"        - per-client sales broken down by month (+ TOTAL row)
#   3) "CUMPLIMIENTO MENSUAL" - per-group rollup (PRESUPUESTO / VENTA / POR CUMPLIR / CUMPLIMIENTO)
#
# Two underlying source values increased (PIEDRA SINTERIZADA sale for
# AUCANSHALA ALLAICA FREDDY HERNAN, PORCELANATO sale for FUENTES PAREDES
# MARIA FERNANDA, both booked in "octubre"), and every dependent total /
# rollup cell is updated to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) VENTAS POR GRUPO - per product-group source values
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("L5").Value  = 427.68     # PIEDRA SINTERIZADA - AUCANSHALA ALLAICA FREDDY HERNAN
$wsGrupo.Range("M13").Value = 1124.02    # PORCELANATO        - FUENTES PAREDES MARIA FERNANDA

# ---------------------------------------------------------------------
# 2) VENTA MENSUAL - "octubre" column for the same two sales, plus the
#    TOTAL row (row 37) which sums the octubre column
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F5").Value  = 476.59    # octubre - AUCANSHALA ALLAICA FREDDY HERNAN
$wsMensual.Range("F13").Value = 1314.1    # octubre - FUENTES PAREDES MARIA FERNANDA
$wsMensual.Range("F37").Value = 44698.85  # octubre TOTAL

# ---------------------------------------------------------------------
# 3) CUMPLIMIENTO MENSUAL - rollups by product group
#    D = VENTA, E = POR CUMPLIR (PRESUPUESTO - VENTA), F = CUMPLIMIENTO (VENTA / PRESUPUESTO)
# ---------------------------------------------------------------------
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 11: PIEDRA SINTERIZADA
$wsCumplimiento.Range("D11").Value = 7577.46
$wsCumplimiento.Range("E11").Value = -4655.23541814726
$wsCumplimiento.Range("F11").Value = 2.593045054461818

# Row 12: PORCELANATO
$wsCumplimiento.Range("D12").Value = 35105.98
$wsCumplimiento.Range("E12").Value = -13404.71
$wsCumplimiento.Range("F12").Value = 1.617692420766158

# Row 14: TOTAL
$wsCumplimiento.Range("D14").Value = 46100.75
$wsCumplimiento.Range("E14").Value = -9515.182762818173
$wsCumplimiento.Range("F14").Value = 1.260080230576497
